$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new values for the "bucket with session" row (row 7): Validation, Exceptions
$ws.Range("F7").Value = "Validation"
$ws.Range("G7").Value = "Exceptions"

# Update the active selection to G8
$ws.Range("G8").Select()
